# issue #5: add legislator_id, name, date into dataframe
#
# Adds three new columns (date, legislator_name, legislator_id) to the
# "股票" (stocks) worksheet, filled in with this filing's metadata:
#   date             = 2011-11-21
#   legislator_name  = 孫大千 (already present elsewhere in the workbook)
#   legislator_id    = 919

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Header row (row 1): bold / centered / bordered, like the existing headers ---
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

$headerRange = $ws.Range($ws.Cells.Item(1, 8), $ws.Cells.Item(1, 10))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data row (row 2): plain style, matching the rest of the row ---
# Force the date to be stored as literal text ("2011-11-21"), not an
# auto-converted date serial number, then drop back to a plain/general
# format so the cell matches its neighbours.
$ws.Cells.Item(2, 8).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value = "2011-11-21"
$ws.Cells.Item(2, 8).ClearFormats()

$ws.Cells.Item(2, 9).Value = "孫大千"
$ws.Cells.Item(2, 10).Value = 919
